# Add a new "添加日志支持" (add logging support) row to the feature list
# (commit: "add log for this app, use log4net").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: new feature entry (#6) ---------------------------------------
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "添加日志支持"

# --- Header row (A1:E1) formatting rotation -------------------------------
# Excel re-saved the header formatting entries in a rotated order even
# though the visible look of each cell is unchanged; reproduce that by
# cycling the per-cell formats: A1:B1 <- E1, C1:D1 <- (old A1:B1), E1 <- (old C1:D1)
$tmp1 = $ws.Range("Z50")
$tmp2 = $ws.Range("Z51")

$ws.Range("A1").Copy()
$tmp1.PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("C1").Copy()
$tmp2.PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("E1").Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)

$tmp1.Copy()
$ws.Range("C1:D1").PasteSpecial(-4122)

$tmp2.Copy()
$ws.Range("E1").PasteSpecial(-4122)

# clean up the temporary holding cells so the used range / dimension is unaffected
$ws.Rows.Item(50).Delete()
$ws.Rows.Item(50).Delete()

# --- Selection -------------------------------------------------------------
$ws.Range("B8").Select()
